# Commit: "29.06.19 Today Sales Details"
#
# The workbook tracks daily cash-denomination counts for several people.
# The "Raju Ahamed" sheet is being updated for the 29.06.19 report:
#   - the "Date: 27.06.19" header (used twice, rows 4 and 31) becomes
#     "Date: 29.06.19"
#   - the Qty. column for the note denominations is refreshed with the
#     day's actual counts (both the top block rows 6-10 and the duplicate
#     block rows 33-37 underneath mirror each other)
#   - the 20-taka row's count is not in yet, so it is cleared back to blank
#   - the sheet's on-screen selection is left on the whole printed range

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Raju Ahamed")

# Update the report date label (appears once above each of the two copies
# of the table on this sheet).
$ws.Range("A4").Value = "Date: 29.06.19"
$ws.Range("A31").Value = "Date: 29.06.19"

# Top table (rows 6-11 correspond to notes 1000/500/100/50/20/10).
$ws.Range("E6").Value = 67
$ws.Range("E7").Value = 308
$ws.Range("E8").Value = 300
$ws.Range("E9").Value = 80
$ws.Range("E10").Value = ""

# Duplicate table further down the sheet (rows 33-38), kept in sync with
# the top one.
$ws.Range("E33").Value = 67
$ws.Range("E34").Value = 308
$ws.Range("E35").Value = 300
$ws.Range("E36").Value = 80
$ws.Range("E37").Value = ""

# Recalculate so the cached SUM() formula results are refreshed.
$excel.Calculate() | Out-Null

# Leave the sheet selected on its full printed range, as last left by the
# author.
$ws.Activate() | Out-Null
$ws.Range("A1:J47").Select() | Out-Null
